$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.653870582580566
$ws.Range("B1").Value = 3.598509311676025
$ws.Range("C1").Value = 1.219458222389221
$ws.Range("D1").Value = 0.9763553142547607
$ws.Range("E1").Value = 0.4960718750953674
